$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adicionando as notas da Lais (linhas 17 e 18)
$ws.Range("C17").Value = 4.25
$ws.Range("C18").Value = 6.2
$ws.Range("F18").Value = 7.4

# Reflete onde o cursor ficou apos o usuario inserir as notas
[void]$ws.Range("C8").Select()
